# Updates cryptos list values (Price / Volume(1h) columns, plus a couple of
# row re-orderings where Coin/Link/Price also changed) to match the latest
# coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "69.081.32"
Set-TextValue "E2" "  -0.99%  "

# Row 3
Set-TextValue "D3" "3.515.44"
Set-TextValue "E3" "  -1.65%  "

# Row 4
Set-TextValue "E4" "  -0.08%  "

# Row 5
Set-TextValue "D5" "572.17"
Set-TextValue "E5" "  -0.86%  "

# Row 6
Set-TextValue "D6" "185.00"
Set-TextValue "E6" "  -2.56%  "

# Row 7
Set-TextValue "D7" "3.506.90"
Set-TextValue "E7" "  -1.84%  "

# Row 8
Set-TextValue "E8" "  -2.78%  "

# Row 9
Set-TextValue "E9" "  +0.03%  "

# Row 10
Set-TextValue "E10" "  +4.07%  "

# Row 11
Set-TextValue "D11" "0.653"
Set-TextValue "E11" "  -1.07%  "

# Row 12
Set-TextValue "D12" "54.26"
Set-TextValue "E12" "  -2.61%  "

# Row 13
Set-TextValue "D13" "0.0000302"
Set-TextValue "E13" "  -0.64%  "

# Row 14
Set-TextValue "D14" "9.46"
Set-TextValue "E14" "  -1.84%  "

# Row 15
Set-TextValue "D15" "4.075.79"
Set-TextValue "E15" "  -1.87%  "

# Row 16
Set-TextValue "D16" "19.45"
Set-TextValue "E16" "  -2.05%  "

# Row 17
Set-TextValue "D17" "3.513.25"
Set-TextValue "E17" "  -1.66%  "

# Row 18
Set-TextValue "D18" "68.948.34"
Set-TextValue "E18" "  -1.22%  "

# Row 19
Set-TextValue "D19" "12.28"
Set-TextValue "E19" "  -2.79%  "

# Row 20
Set-TextValue "E20" "  -1.10%  "

# Row 21
Set-TextValue "D21" "546.59"
Set-TextValue "E21" "  +14.58%  "

# Row 22
Set-TextValue "E22" "  -2.77%  "

# Row 23
Set-TextValue "D23" "18.69"
Set-TextValue "E23" "  -4.09%  "

# Row 24
Set-TextValue "E24" "  -1.11%  "

# Row 25
Set-TextValue "D25" "4.41"
Set-TextValue "E25" "  +0.96%  "

# Row 26
Set-TextValue "D26" "94.11"
Set-TextValue "E26" "  -1.73%  "

# Row 27
Set-TextValue "E27" "  -2.56%  "

# Row 28
Set-TextValue "D28" "10.83"
Set-TextValue "E28" "  -1.83%  "

# Row 29
Set-TextValue "D29" "9.12"
Set-TextValue "E29" "  -2.69%  "

# Row 30
Set-TextValue "D30" "31.89"
Set-TextValue "E30" "  -1.32%  "

# Row 31
Set-TextValue "D31" "7.25"
Set-TextValue "E31" "  -5.36%  "

# Row 32
Set-TextValue "D32" "12.59"
Set-TextValue "E32" "  +3.07%  "

# Row 33
Set-TextValue "D33" "64.92"
Set-TextValue "E33" "  -2.07%  "

# Row 34
Set-TextValue "E34" "  -3.76%  "

# Row 35
Set-TextValue "D35" "561.52"
Set-TextValue "E35" "  -3.60%  "

# Row 36
Set-TextValue "D36" "37.98"
Set-TextValue "E36" "  -2.34%  "

# Row 37
Set-TextValue "E37" "  +0.02%  "

# Row 38
Set-TextValue "D38" "0.399"
Set-TextValue "E38" "  +0.92%  "

# Row 39
Set-TextValue "D39" "3.04"
Set-TextValue "E39" "  +6.42%  "

# Row 40
Set-TextValue "D40" "0.0₃0765"
Set-TextValue "E40" "  -4.45%  "

# Row 41
Set-TextValue "D41" "3.14"
Set-TextValue "E41" "  -2.93%  "

# Row 42
Set-TextValue "B42" "Stacks"
Set-TextValue "C42" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D42" "3.36"
Set-TextValue "E42" "  -2.90%  "

# Row 43
Set-TextValue "B43" "Kaspa"
Set-TextValue "C43" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D43" "0.133"
Set-TextValue "E43" "  -3.29%  "

# Row 44
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "3.264.80"
Set-TextValue "E44" "  +1.15%  "

# Row 45
Set-TextValue "B45" "ApeXProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D45" "3.53"
Set-TextValue "E45" "  +5.00%  "

# Row 46
Set-TextValue "E46" "  -2.75%  "

# Row 47
Set-TextValue "D47" "0.0441"
Set-TextValue "E47" "  -0.35%  "

# Row 48
Set-TextValue "D48" "0.135"
Set-TextValue "E48" "  -2.32%  "

# Row 49
Set-TextValue "D49" "8.96"
Set-TextValue "E49" "  -3.74%  "

# Row 50
Set-TextValue "D50" "0.998"

# Row 51
Set-TextValue "D51" "137.57"
Set-TextValue "E51" "  +3.10%  "
